$d = $word.ActiveDocument

# Locate the placeholder "R$ 0,00" value inside the "Custo Estimado" cell.
$rng = $d.Content
$rng.Find.Execute("R$ 0,00")

if ($rng.Find.Found) {
    # Keep the "R$ " prefix, collapse to the end of it, then insert the
    # actual amount as its own run (mirrors how Word splits a run when
    # new text is typed in place of the old value).
    $rng.Text = "R$ "
    $rng.Collapse(0)
    $rng.InsertAfter("12.483,20")

    # Force the newly inserted text onto its own run (distinct w:r) while
    # leaving its formatting identical to the "R$ " run.
    $rng.Bold = 1
    $rng.Bold = 0
}
